# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the 30c1937e... file moves
# from "Handed back: in sync with en-US" to "In Translation", the
# a3f10386... file moves to "Ready for handoff", handoff/handback
# timestamps are refreshed, and an "Error Detail" / "Latest Handback
# DateTime" style note is recorded on the per-locale sheets explaining
# that the last handback isn't against the latest source.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errMsg30c = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77045efd0b02209e77186843a8991ea79a43c299/e2e/30c1937e-af3c-4537-8f26-9b07f24af10f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f98598600d282654e0a8446fca6475dae4e0f422/e2e/30c1937e-af3c-4537-8f26-9b07f24af10f.md."
$errMsgA3f = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77045efd0b02209e77186843a8991ea79a43c299/e2e/a3f10386-b88f-4224-b418-e3ad96775d41.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f98598600d282654e0a8446fca6475dae4e0f422/e2e/a3f10386-b88f-4224-b418-e3ad96775d41.md."

# --- Overview sheet -------------------------------------------------------
# Columns: A File Name | B Path And Name | C Extension | D Publish URL |
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
# Row 2 = 30c1937e-af3c-4537-8f26-9b07f24af10f.md, Row 3 = a3f10386-...md

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("G2").Value = "2016-09-07 03:40:04"

$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 03:40:04"

# --- zh-cn sheet ------------------------------------------------------------
# Columns: A Source File Name | B File Extension | C Status | D Source Path |
#          E Priority | F Content Duplicate | G Latest Handoff File |
#          H Latest Handoff Datetime | I Latest Target File |
#          J Latest Handback File | K Latest Handback DateTime |
#          L Reference Tokens | M To be localized | N Dependency From |
#          O Has metadata | P Error Detail

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("H2").Value = "2016-09-07 03:39:51"
$zhcn.Range("P2").Value = $errMsg30c

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-07 03:39:51"
$zhcn.Range("P3").Value = $errMsgA3f

# --- de-de sheet ------------------------------------------------------------
# Same column layout as zh-cn.

$dede.Range("C2").Value = "In Translation"
$dede.Range("H2").Value = "2016-09-07 03:40:04"
$dede.Range("P2").Value = $errMsg30c

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-07 03:40:04"
$dede.Range("P3").Value = $errMsgA3f

# --- Column width touch-ups -------------------------------------------------
# Narrower zh-cn/de-de status columns on Overview and the per-locale Status
# column, and a wider Error Detail column to host the new long message.
# (ColumnWidth is quantized to 1/6-character steps by this host, so we land
# on the closest reachable step to the authored width.)

$overview.Columns.Item(5).ColumnWidth = 16.3   # E -> ~17.22
$overview.Columns.Item(6).ColumnWidth = 16.3   # F -> ~17.22

$zhcn.Columns.Item(3).ColumnWidth = 16.3        # C -> ~17.22
$zhcn.Columns.Item(16).ColumnWidth = 39.1667    # P -> 40

$dede.Columns.Item(3).ColumnWidth = 16.3        # C -> ~17.22
$dede.Columns.Item(16).ColumnWidth = 39.1667    # P -> 40
